$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "D" column (Price) cells as plain text, preserving exact formatting ---
# Force text storage (avoid Excel auto-converting numeric-looking strings to numbers,
# which would drop significant trailing zeros / dot-grouping), then restore default style.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.032.71"
$ws.Range("D3").Value = "3.190.26"
$ws.Range("D5").Value = "604.46"
$ws.Range("D6").Value = "153.70"
$ws.Range("D7").Value = "0.999"
$ws.Range("D8").Value = "3.191.11"
$ws.Range("D11").Value = "5.62"
$ws.Range("D12").Value = "0.482"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D14").Value = "37.51"
$ws.Range("D15").Value = "3.682.83"
$ws.Range("D16").Value = "65.044.10"
$ws.Range("D17").Value = "3.186.36"
$ws.Range("D18").Value = "0.114"
$ws.Range("D19").Value = "7.12"
$ws.Range("D20").Value = "486.73"
$ws.Range("D21").Value = "15.01"
$ws.Range("D22").Value = "0.724"
$ws.Range("D23").Value = "7.88"
$ws.Range("D24").Value = "14.18"
$ws.Range("D25").Value = "85.12"
$ws.Range("D27").Value = "2.96"
$ws.Range("D28").Value = "8.87"
$ws.Range("D29").Value = "2.29"
$ws.Range("D30").Value = "7.28"
$ws.Range("D31").Value = "0.123"
$ws.Range("D32").Value = "2.75"
$ws.Range("D33").Value = "27.13"
$ws.Range("D36").Value = "6.20"
$ws.Range("D37").Value = "3.35"
$ws.Range("D38").Value = "54.66"
$ws.Range("D39").Value = "0.0₃0756"
$ws.Range("D40").Value = "466.10"
$ws.Range("D41").Value = "0.129"
$ws.Range("D42").Value = "0.0409"
$ws.Range("D43").Value = "8.60"
$ws.Range("D44").Value = "2.49"
$ws.Range("D45").Value = "2.944.90"
$ws.Range("D46").Value = "0.280"
$ws.Range("D47").Value = "27.62"

foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}

# --- Update "B"/"C"/"E" columns (Coin name, Link, Volume) - plain text values ---
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  -3.98%  "
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("E12").Value = "  -4.90%  "
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  -3.82%  "
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("E31").Value = "  +2.17%  "
$ws.Range("E32").Value = "  -7.11%  "
$ws.Range("E33").Value = "  -3.77%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  -4.95%  "
$ws.Range("E36").Value = "  -5.00%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E37").Value = "  +9.38%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("E40").Value = "  -8.32%  "
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("E46").Value = "  -6.66%  "
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -1.90%  "
